# Update "想去人数" (interested count) figures for the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2328
    "F3" = 1784
    "F4" = 345
    "F6" = 975
    "F7" = 43
    "F8" = 5894
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
